$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.237540118649887 ; $ws.Range("C2").Value = 4.847265884137244 ; $ws.Range("E2").Value = 20.45706822967662 ; $ws.Range("F2").Value = 40.60758663165755 ; $ws.Range("G2").Value = 3.648126559369127 ; $ws.Range("J2").Value = 8.29605291969596 ; $ws.Range("K2").Value = 8.566395636622991 ; $ws.Range("M2").Value = 17.79634646392574 ; $ws.Range("N2").Value = 19.36435176894801 ; $ws.Range("O2").Value = 23.4633967475671
$ws.Range("B3").Value = 8.995049700448309 ; $ws.Range("C3").Value = 4.739895602674514 ; $ws.Range("E3").Value = 20.3522692708156 ; $ws.Range("F3").Value = 40.54750959811948 ; $ws.Range("G3").Value = 3.649765675249533 ; $ws.Range("J3").Value = 8.309655889786097 ; $ws.Range("K3").Value = 8.409719272221309 ; $ws.Range("M3").Value = 17.68671362961611 ; $ws.Range("N3").Value = 19.4251663442612 ; $ws.Range("O3").Value = 23.54381649680408
$ws.Range("B4").Value = 8.844512339022415 ; $ws.Range("C4").Value = 4.672093838260674 ; $ws.Range("E4").Value = 20.29193732667622 ; $ws.Range("F4").Value = 40.52049975945128 ; $ws.Range("G4").Value = 3.650825330958281 ; $ws.Range("J4").Value = 8.318432438427561 ; $ws.Range("K4").Value = 8.313419467581507 ; $ws.Range("M4").Value = 17.6222911483426 ; $ws.Range("N4").Value = 19.46424754347827 ; $ws.Range("O4").Value = 23.59776008380815
$ws.Range("B5").Value = 8.782848998434305 ; $ws.Range("C5").Value = 4.644015880806089 ; $ws.Range("E5").Value = 20.26838022869093 ; $ws.Range("F5").Value = 40.51198469140743 ; $ws.Range("G5").Value = 3.651270576309386 ; $ws.Range("J5").Value = 8.322115948275222 ; $ws.Range("K5").Value = 8.27420317862633 ; $ws.Range("M5").Value = 17.59678741466305 ; $ws.Range("N5").Value = 19.48061250271804 ; $ws.Range("O5").Value = 23.62088885526716
$ws.Range("B6").Value = 8.772593533881185 ; $ws.Range("C6").Value = 4.639327202487727 ; $ws.Range("E6").Value = 20.26453124653757 ; $ws.Range("F6").Value = 40.51072144756427 ; $ws.Range("G6").Value = 3.651345321140878 ; $ws.Range("J6").Value = 8.322734063601875 ; $ws.Range("K6").Value = 8.267694484084027 ; $ws.Range("M6").Value = 17.59259838274977 ; $ws.Range("N6").Value = 19.48335644812758 ; $ws.Range("O6").Value = 23.62479855431643
$ws.Range("B7").Value = 8.843681883326184 ; $ws.Range("C7").Value = 4.67171695169772 ; $ws.Range("E7").Value = 20.29161543836372 ; $ws.Range("F7").Value = 40.52037482492175 ; $ws.Range("G7").Value = 3.650831281272084 ; $ws.Range("J7").Value = 8.318481681900877 ; $ws.Range("K7").Value = 8.312890405990398 ; $ws.Range("M7").Value = 17.62194413591911 ; $ws.Range("N7").Value = 19.46446646759409 ; $ws.Range("O7").Value = 23.59806736740583
$ws.Range("B8").Value = 9.154327424779032 ; $ws.Range("C8").Value = 4.810647666851724 ; $ws.Range("E8").Value = 20.42011209125424 ; $ws.Range("F8").Value = 40.58482728848302 ; $ws.Range("G8").Value = 3.648680703230974 ; $ws.Range("J8").Value = 8.300655353102911 ; $ws.Range("K8").Value = 8.51242512635657 ; $ws.Range("M8").Value = 17.75795801646222 ; $ws.Range("N8").Value = 19.38496012513485 ; $ws.Range("O8").Value = 23.49017687743467
$ws.Range("B9").Value = 9.746471210555852 ; $ws.Range("C9").Value = 5.067242695378049 ; $ws.Range("E9").Value = 20.7030487063864 ; $ws.Range("F9").Value = 40.78915649016181 ; $ws.Range("G9").Value = 3.644883903759943 ; $ws.Range("J9").Value = 8.269050145759211 ; $ws.Range("K9").Value = 8.900699850601102 ; $ws.Range("M9").Value = 18.04663265360691 ; $ws.Range("N9").Value = 19.24279995322832 ; $ws.Range("O9").Value = 23.31490528343776
$ws.Range("B10").Value = 10.16622938056749 ; $ws.Range("C10").Value = 5.244914831904653 ; $ws.Range("E10").Value = 20.92849589262885 ; $ws.Range("F10").Value = 40.9860699686698 ; $ws.Range("G10").Value = 3.642348074480376 ; $ws.Range("J10").Value = 8.247853280834148 ; $ws.Range("K10").Value = 9.181259742303009 ; $ws.Range("M10").Value = 18.27070147473968 ; $ws.Range("N10").Value = 19.14665136749751 ; $ws.Range("O10").Value = 23.20836568086163
$ws.Range("B11").Value = 10.35294976246302 ; $ws.Range("C11").Value = 5.323164429085383 ; $ws.Range("E11").Value = 21.03457246938648 ; $ws.Range("F11").Value = 41.08561324094723 ; $ws.Range("G11").Value = 3.641248979367338 ; $ws.Range("J11").Value = 8.238645375004189 ; $ws.Range("K11").Value = 9.30730118980577 ; $ws.Range("M11").Value = 18.37492186756184 ; $ws.Range("N11").Value = 19.10469385964631 ; $ws.Range("O11").Value = 23.16474685641702
$ws.Range("B12").Value = 10.42297943166922 ; $ws.Range("C12").Value = 5.352408593078978 ; $ws.Range("E12").Value = 21.0752200671565 ; $ws.Range("F12").Value = 41.12472060398306 ; $ws.Range("G12").Value = 3.640840570169058 ; $ws.Range("J12").Value = 8.235220763397214 ; $ws.Range("K12").Value = 9.35475701355449 ; $ws.Range("M12").Value = 18.41468923697994 ; $ws.Range("N12").Value = 19.08906037517744 ; $ws.Range("O12").Value = 23.1489281982012
$ws.Range("B13").Value = 10.40792844718176 ; $ws.Range("C13").Value = 5.346127819731898 ; $ws.Range("E13").Value = 21.06644503464265 ; $ws.Range("F13").Value = 41.116235663802 ; $ws.Range("G13").Value = 3.640928182345438 ; $ws.Range("J13").Value = 8.23595555185728 ; $ws.Range("K13").Value = 9.344549402922102 ; $ws.Range("M13").Value = 18.40611166929495 ; $ws.Range("N13").Value = 19.09241600723921 ; $ws.Range("O13").Value = 23.15230391812278
$ws.Range("B14").Value = 10.35872509336896 ; $ws.Range("C14").Value = 5.325578218190218 ; $ws.Range("E14").Value = 21.03790713781746 ; $ws.Range("F14").Value = 41.08880244458539 ; $ws.Range("G14").Value = 3.641215223303747 ; $ws.Range("J14").Value = 8.238362384711859 ; $ws.Range("K14").Value = 9.311211137172034 ; $ws.Range("M14").Value = 18.37818767498527 ; $ws.Range("N14").Value = 19.10340258139513 ; $ws.Range("O14").Value = 23.16343142924829
$ws.Range("B15").Value = 10.32849646958573 ; $ws.Range("C15").Value = 5.312940078180752 ; $ws.Range("E15").Value = 21.0204883650907 ; $ws.Range("F15").Value = 41.07218210465734 ; $ws.Range("G15").Value = 3.641392058102009 ; $ws.Range("J15").Value = 8.239844734257968 ; $ws.Range("K15").Value = 9.290753590960419 ; $ws.Range("M15").Value = 18.36112183098749 ; $ws.Range("N15").Value = 19.11016533745261 ; $ws.Range("O15").Value = 23.17033841707241
$ws.Range("B16").Value = 10.15393560960962 ; $ws.Range("C16").Value = 5.239747743080292 ; $ws.Range("E16").Value = 20.9216320381059 ; $ws.Range("F16").Value = 40.97976341117947 ; $ws.Range("G16").Value = 3.642420995684138 ; $ws.Range("J16").Value = 8.248463761480146 ; $ws.Range("K16").Value = 9.17298679261155 ; $ws.Range("M16").Value = 18.26393413371887 ; $ws.Range("N16").Value = 19.14942912513204 ; $ws.Range("O16").Value = 23.21131396377223
$ws.Range("B17").Value = 10.04571377444609 ; $ws.Range("C17").Value = 5.194174673782543 ; $ws.Range("E17").Value = 20.8618691024369 ; $ws.Range("F17").Value = 40.92560611653553 ; $ws.Range("G17").Value = 3.643066138616422 ; $ws.Range("J17").Value = 8.25386238076107 ; $ws.Range("K17").Value = 9.100301389139446 ; $ws.Range("M17").Value = 18.20487949227918 ; $ws.Range("N17").Value = 19.17397148680582 ; $ws.Range("O17").Value = 23.2376938358014
$ws.Range("B18").Value = 9.98307305486645 ; $ws.Range("C18").Value = 5.167720870731246 ; $ws.Range("E18").Value = 20.82782801623571 ; $ws.Range("F18").Value = 40.89539572300162 ; $ws.Range("G18").Value = 3.643442336644378 ; $ws.Range("J18").Value = 8.257008452113025 ; $ws.Range("K18").Value = 9.058347877137088 ; $ws.Range("M18").Value = 18.17113011326204 ; $ws.Range("N18").Value = 19.18825530734286 ; $ws.Range("O18").Value = 23.25332291182958
$ws.Range("B19").Value = 9.961798451969681 ; $ws.Range("C19").Value = 5.158723153504609 ; $ws.Range("E19").Value = 20.81636029843585 ; $ws.Range("F19").Value = 40.88532894979955 ; $ws.Range("G19").Value = 3.643570592883695 ; $ws.Range("J19").Value = 8.258080695875401 ; $ws.Range("K19").Value = 9.044119343507699 ; $ws.Range("M19").Value = 18.15974131769355 ; $ws.Range("N19").Value = 19.19312040311563 ; $ws.Range("O19").Value = 23.2586929200271
$ws.Range("B20").Value = 10.05727551015902 ; $ws.Range("C20").Value = 5.199051114700715 ; $ws.Range("E20").Value = 20.86819670316532 ; $ws.Range("F20").Value = 40.93127416651866 ; $ws.Range("G20").Value = 3.642996931505183 ; $ws.Range("J20").Value = 8.253283454630321 ; $ws.Range("K20").Value = 9.108054401616556 ; $ws.Range("M20").Value = 18.21114367229638 ; $ws.Range("N20").Value = 19.17134156248049 ; $ws.Range("O20").Value = 23.23483843767866
$ws.Range("B21").Value = 10.3731962128011 ; $ws.Range("C21").Value = 5.33162477360996 ; $ws.Range("E21").Value = 21.04627664067427 ; $ws.Range("F21").Value = 41.09682208646585 ; $ws.Range("G21").Value = 3.641130701183433 ; $ws.Range("J21").Value = 8.237653753081599 ; $ws.Range("K21").Value = 9.321011160194276 ; $ws.Range("M21").Value = 18.38638167518522 ; $ws.Range("N21").Value = 19.10016864938729 ; $ws.Range("O21").Value = 23.16014402809419
$ws.Range("B22").Value = 10.57569311212897 ; $ws.Range("C22").Value = 5.416005123743232 ; $ws.Range("E22").Value = 21.1654391268912 ; $ws.Range("F22").Value = 41.21323925115052 ; $ws.Range("G22").Value = 3.63995642536384 ; $ws.Range("J22").Value = 8.227801391485229 ; $ws.Range("K22").Value = 9.458578281063252 ; $ws.Range("M22").Value = 18.50265317687375 ; $ws.Range("N22").Value = 19.0551383218971 ; $ws.Range("O22").Value = 23.11540115969488
$ws.Range("B23").Value = 10.4680017261487 ; $ws.Range("C23").Value = 5.371182201357833 ; $ws.Range("E23").Value = 21.10159503901397 ; $ws.Range("F23").Value = 41.1503602576351 ; $ws.Range("G23").Value = 3.640579015573073 ; $ws.Range("J23").Value = 8.233026700772838 ; $ws.Range("K23").Value = 9.385317947882095 ; $ws.Range("M23").Value = 18.44044672736674 ; $ws.Range("N23").Value = 19.07903634099861 ; $ws.Range("O23").Value = 23.13890785897459
$ws.Range("B24").Value = 10.05204976025573 ; $ws.Range("C24").Value = 5.196847263256571 ; $ws.Range("E24").Value = 20.8653350040509 ; $ws.Range("F24").Value = 40.92870875712454 ; $ws.Range("G24").Value = 3.643028203538369 ; $ws.Range("J24").Value = 8.253545055246779 ; $ws.Range("K24").Value = 9.104549777338518 ; $ws.Range("M24").Value = 18.208311005408 ; $ws.Range("N24").Value = 19.17253000888165 ; $ws.Range("O24").Value = 23.23612792127577
$ws.Range("B25").Value = 9.588633597509727 ; $ws.Range("C25").Value = 4.99965168562037 ; $ws.Range("E25").Value = 20.62332187445142 ; $ws.Range("F25").Value = 40.72560190509649 ; $ws.Range("G25").Value = 3.645866299214422 ; $ws.Range("J25").Value = 8.277243404303626 ; $ws.Range("K25").Value = 8.796286811971081 ; $ws.Range("M25").Value = 18.04663265360691 ; $ws.Range("N25").Value = 19.2797948643526 ; $ws.Range("O25").Value = 23.35842388314468
